$d = $word.ActiveDocument
$nl = [char]11

# Replace the 15 lattice-multiplication exercises (5 rows x 3 cols) in
# the single table with a new set of problems/answers. Each cell's
# paragraph run holds 5 lines (problem, top factors, separator, two
# partial-product placeholders) joined by <w:br/> line breaks, which we
# reproduce here as vertical-tab (chr 11) separated text so Word
# round-trips them back to <w:br/> elements. Whole-cell replacement
# (rather than piecemeal Find/Replace) avoids cross-line text
# collisions where a new value for one line matches another line's
# still-unreplaced old value.

# Cell row 1, col 1: '71 x 11' -> '52 x 61'
$cell = $d.Tables(1).Cell(1, 1)
$cell.Range.Text = "52 x 61" + $nl + "  6    1" + $nl + "  ----" + $nl + "5|    |" + $nl + "2|    |"

# Cell row 1, col 2: '49 x 81' -> '99 x 46'
$cell = $d.Tables(1).Cell(1, 2)
$cell.Range.Text = "99 x 46" + $nl + "  4    6" + $nl + "  ----" + $nl + "9|    |" + $nl + "9|    |"

# Cell row 1, col 3: '83 x 57' -> '48 x 14'
$cell = $d.Tables(1).Cell(1, 3)
$cell.Range.Text = "48 x 14" + $nl + "  1    4" + $nl + "  ----" + $nl + "4|    |" + $nl + "8|    |"

# Cell row 2, col 1: '98 x 56' -> '56 x 40'
$cell = $d.Tables(1).Cell(2, 1)
$cell.Range.Text = "56 x 40" + $nl + "  4    0" + $nl + "  ----" + $nl + "5|    |" + $nl + "6|    |"

# Cell row 2, col 2: '95 x 66' -> '28 x 32'
$cell = $d.Tables(1).Cell(2, 2)
$cell.Range.Text = "28 x 32" + $nl + "  3    2" + $nl + "  ----" + $nl + "2|    |" + $nl + "8|    |"

# Cell row 2, col 3: '31 x 47' -> '37 x 67'
$cell = $d.Tables(1).Cell(2, 3)
$cell.Range.Text = "37 x 67" + $nl + "  6    7" + $nl + "  ----" + $nl + "3|    |" + $nl + "7|    |"

# Cell row 3, col 1: '14 x 29' -> '52 x 63'
$cell = $d.Tables(1).Cell(3, 1)
$cell.Range.Text = "52 x 63" + $nl + "  6    3" + $nl + "  ----" + $nl + "5|    |" + $nl + "2|    |"

# Cell row 3, col 2: '73 x 66' -> '34 x 29'
$cell = $d.Tables(1).Cell(3, 2)
$cell.Range.Text = "34 x 29" + $nl + "  2    9" + $nl + "  ----" + $nl + "3|    |" + $nl + "4|    |"

# Cell row 3, col 3: '71 x 61' -> '12 x 76'
$cell = $d.Tables(1).Cell(3, 3)
$cell.Range.Text = "12 x 76" + $nl + "  7    6" + $nl + "  ----" + $nl + "1|    |" + $nl + "2|    |"

# Cell row 4, col 1: '30 x 69' -> '82 x 81'
$cell = $d.Tables(1).Cell(4, 1)
$cell.Range.Text = "82 x 81" + $nl + "  8    1" + $nl + "  ----" + $nl + "8|    |" + $nl + "2|    |"

# Cell row 4, col 2: '63 x 91' -> '25 x 43'
$cell = $d.Tables(1).Cell(4, 2)
$cell.Range.Text = "25 x 43" + $nl + "  4    3" + $nl + "  ----" + $nl + "2|    |" + $nl + "5|    |"

# Cell row 4, col 3: '31 x 24' -> '30 x 42'
$cell = $d.Tables(1).Cell(4, 3)
$cell.Range.Text = "30 x 42" + $nl + "  4    2" + $nl + "  ----" + $nl + "3|    |" + $nl + "0|    |"

# Cell row 5, col 1: '76 x 72' -> '27 x 80'
$cell = $d.Tables(1).Cell(5, 1)
$cell.Range.Text = "27 x 80" + $nl + "  8    0" + $nl + "  ----" + $nl + "2|    |" + $nl + "7|    |"

# Cell row 5, col 2: '41 x 66' -> '13 x 37'
$cell = $d.Tables(1).Cell(5, 2)
$cell.Range.Text = "13 x 37" + $nl + "  3    7" + $nl + "  ----" + $nl + "1|    |" + $nl + "3|    |"

# Cell row 5, col 3: '21 x 83' -> '84 x 39'
$cell = $d.Tables(1).Cell(5, 3)
$cell.Range.Text = "84 x 39" + $nl + "  3    9" + $nl + "  ----" + $nl + "8|    |" + $nl + "4|    |"

